$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- E3: Quantity value changed from 5 to 14 ---
$ws.Range("E3").Value = 14

# --- Add new row 4 (mirrors rows 2/3: a book entry) ---
$ws.Range("A4").Value = "intrstelr.007@gmail.com"
$ws.Range("B4").Value = "anujsharma.cv"
$ws.Range("C4").Value = "asfsdsfs"
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "NS10"

# Add mailto hyperlink on A4, then restyle to match the existing
# hyperlink cells (A2/A3) so it shares their "Hyperlink" cell style.
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:intrstelr.007@gmail.com")
$ws.Range("A4").Style = $ws.Range("A3").Style

# --- Selection moves to E5 (next empty row) ---
[void]$ws.Range("E5").Select()
